$d = $word.ActiveDocument

$replacements = @(
    @("429×8=3432", "837×8=6696"),
    @("428×3=1284", "419×8=3352"),
    @("255×5=1275", "248×4=992"),
    @("480×7=3360", "676×8=5408"),
    @("635×8=5080", "365×5=1825"),
    @("429×7=3003", "207×2=414"),
    @("558×9=5022", "115×4=460"),
    @("651×6=3906", "389×5=1945"),
    @("958×2=1916", "780×9=7020"),
    @("860×8=6880", "457×4=1828"),
    @("698×4=2792", "128×2=256"),
    @("566×2=1132", "492×2=984"),
    @("234×8=1872", "616×4=2464"),
    @("624×3=1872", "177×8=1416"),
    @("857×5=4285", "376×5=1880"),
    @("728×8=5824", "665×5=3325"),
    @("508×7=3556", "210×6=1260"),
    @("639×5=3195", "919×8=7352"),
    @("733×3=2199", "570×8=4560"),
    @("505×4=2020", "674×4=2696"),
    @("108×2=216",  "841×7=5887"),
    @("577×6=3462", "434×5=2170"),
    @("900×5=4500", "668×3=2004"),
    @("717×2=1434", "820×6=4920"),
    @("238×7=1666", "746×8=5968")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
